$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B23").Value = "   health >= INJURY_CLEAR"
$ws.Range("B24").Value = "   then setInjured(false boolean, playerId int)"
$ws.Range("B25").Value = "   checkHealthMax(health int,playerID int)"
$ws.Range("B27").Value = "   print ""Player is still injured"""
$ws.Range("B28").Value = "   setHealth(health int, playerID)"

$ws.Range("B28").Select()
